# code added if incorrect ruleName in file
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Simulate an incorrect RuleName being supplied for the first test row.
$ws.Range("C2").Value = "xyz"

# Clear the previously-computed Actual / Pass-Fail result columns for both
# data rows so they no longer hold stale results.
$ws.Range("R2").Value = $null
$ws.Range("S2").Value = $null
$ws.Range("R5").Value = $null
$ws.Range("S5").Value = $null

# Update the visible window / selection to match the reviewed range.
$ws.Range("A1").Select()
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("R2:S5").Select()
